$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Concent"
$ws.Range("B12").Value = "Kaymary"
$ws.Range("C12").Value = 500.0
$ws.Range("D12").Value = 23.0
$ws.Range("E12").Value = 1.0
